# Generate Report for Handoff
# Updates the localization-status report: flips the "In Translation" status
# to "Ready for handoff" and refreshes the associated timestamps on all
# three sheets (Overview, zh-cn, de-de), then widens the Status columns to
# fit the new (longer) text, mirroring an auto-generated report refresh.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps refreshed alongside the status change ---
$wsOverview.Range("G2").Value = "2016-08-24 04:37:21"
$wsZhCn.Range("H2").Value = "2016-08-24 04:37:17"
$wsDeDe.Range("H2").Value = "2016-08-24 04:37:21"

# --- Widen the Status columns so the new text fits ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
